$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows above row 214, shifting the existing rows 214-219 down to 216-221.
$ws.Rows.Item(214).Insert()
$ws.Rows.Item(214).Insert()

# Fill in the new row 214 (Americana (o), Primera) with this week's data.
$ws.Cells.Item(214, 1).Value = 2
$ws.Cells.Item(214, 2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(214, 3).Value = "Coquimbo"
$ws.Cells.Item(214, 4).Value = 44595
$ws.Cells.Item(214, 5).Value = 4
$ws.Cells.Item(214, 6).Value = 100112021
$ws.Cells.Item(214, 7).Value = "Ají"
$ws.Cells.Item(214, 8).Value = "Americana (o)"
$ws.Cells.Item(214, 9).Value = "Primera"
$ws.Cells.Item(214, 10).Value = 200
$ws.Cells.Item(214, 11).Value = 10000
$ws.Cells.Item(214, 12).Value = 12000
$ws.Cells.Item(214, 13).Value = 11000
$ws.Cells.Item(214, 14).Value = "$/caja 25 kilos"
$ws.Cells.Item(214, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(214, 16).Value = 440
$ws.Cells.Item(214, 17).Value = 25
$ws.Cells.Item(214, 18).Value = "Hortaliza"

# Fill in the new row 215 (Inferno, Primera) with this week's data.
$ws.Cells.Item(215, 1).Value = 2
$ws.Cells.Item(215, 2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(215, 3).Value = "Coquimbo"
$ws.Cells.Item(215, 4).Value = 44595
$ws.Cells.Item(215, 5).Value = 4
$ws.Cells.Item(215, 6).Value = 100112021
$ws.Cells.Item(215, 7).Value = "Ají"
$ws.Cells.Item(215, 8).Value = "Inferno"
$ws.Cells.Item(215, 9).Value = "Primera"
$ws.Cells.Item(215, 10).Value = 200
$ws.Cells.Item(215, 11).Value = 13000
$ws.Cells.Item(215, 12).Value = 15000
$ws.Cells.Item(215, 13).Value = 14000
$ws.Cells.Item(215, 14).Value = "$/caja 25 kilos"
$ws.Cells.Item(215, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(215, 16).Value = 560
$ws.Cells.Item(215, 17).Value = 25
$ws.Cells.Item(215, 18).Value = "Hortaliza"

# Apply the same date-number format used by the other cells in column D.
$ws.Cells.Item(214, 4).NumberFormat = $ws.Cells.Item(216, 4).NumberFormat
$ws.Cells.Item(215, 4).NumberFormat = $ws.Cells.Item(216, 4).NumberFormat
